# Refresh the cryptos price/volume snapshot (GitHub Actions scheduled update).
# Column layout: A=rank(unchanged) B=Coin C=Link D=Price E=Volume(1h)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Cells.Item(2, 4).Value = '68.849.15'
$ws.Cells.Item(2, 5).Value = '  +0.26%  '

# Row 3: Ethereum
$ws.Cells.Item(3, 4).Value = '3.759.31'
$ws.Cells.Item(3, 5).Value = '  -1.45%  '

# Row 4: TetherUSD
$ws.Cells.Item(4, 5).Value = '  +0.58%  '

# Row 5: BNB
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '631.16'
$ws.Cells.Item(5, 5).Value = '  +2.95%  '

# Row 6: Solana
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '165.10'
$ws.Cells.Item(6, 5).Value = '  +0.37%  '

# Row 7: LidoStakedEther
$ws.Cells.Item(7, 4).Value = '3.753.41'
$ws.Cells.Item(7, 5).Value = '  -1.51%  '

# Row 8: USDC
$ws.Cells.Item(8, 5).Value = '  +0.09%  '

# Row 9: XRP
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.519'
$ws.Cells.Item(9, 5).Value = '  +0.36%  '

# Row 10: Dogecoin
$ws.Cells.Item(10, 5).Value = '  -1.67%  '

# Row 11: Cardano
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.456'
$ws.Cells.Item(11, 5).Value = '  +1.26%  '

# Row 12: Toncoin
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '6.78'
$ws.Cells.Item(12, 5).Value = '  +1.25%  '

# Row 13: ShibaInu
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '0.0000238'
$ws.Cells.Item(13, 5).Value = '  -3.92%  '

# Row 14: Avalanche
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '34.72'
$ws.Cells.Item(14, 5).Value = '  -2.03%  '

# Row 15: WrappedliquidstakedEther2.0
$ws.Cells.Item(15, 4).Value = '4.395.28'
$ws.Cells.Item(15, 5).Value = '  -1.34%  '

# Row 16: WrappedEther
$ws.Cells.Item(16, 4).Value = '3.762.17'
$ws.Cells.Item(16, 5).Value = '  -0.30%  '

# Row 17: WrappedBTC
$ws.Cells.Item(17, 4).Value = '68.912.38'
$ws.Cells.Item(17, 5).Value = '  +0.37%  '

# Row 18: Chainlink
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '17.60'
$ws.Cells.Item(18, 5).Value = '  -2.60%  '

# Row 19: TRON
$ws.Cells.Item(19, 5).Value = '  +0.20%  '

# Row 20: Polkadot
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '7.00'
$ws.Cells.Item(20, 5).Value = '  -1.54%  '

# Row 21: BitcoinCash
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '461.57'
$ws.Cells.Item(21, 5).Value = '  -0.56%  '

# Row 22: Uniswap
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '9.47'
$ws.Cells.Item(22, 5).Value = '  -2.16%  '

# Row 23: Polygon
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '0.703'
$ws.Cells.Item(23, 5).Value = '  +0.23%  '

# Row 24: PEPE -> Litecoin
$ws.Cells.Item(24, 2).Value = 'Litecoin'
$ws.Cells.Item(24, 3).Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '81.97'
$ws.Cells.Item(24, 5).Value = '  -2.17%  '

# Row 25: Litecoin -> PEPE
$ws.Cells.Item(25, 2).Value = 'PEPE'
$ws.Cells.Item(25, 3).Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '0.0000143'
$ws.Cells.Item(25, 5).Value = '  -5.01%  '

# Row 26: InternetComputer(DFINITY)
$ws.Cells.Item(26, 5).Value = '  +0.80%  '

# Row 27: Fetch.AI
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '2.11'
$ws.Cells.Item(27, 5).Value = '  -0.16%  '

# Row 28: RenderToken
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '10.07'
$ws.Cells.Item(28, 5).Value = '  +0.92%  '

# Row 30: WrappedeETH
$ws.Cells.Item(30, 4).Value = '3.911.38'
$ws.Cells.Item(30, 5).Value = '  -1.35%  '

# Row 31: ImmutableX
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '2.28'
$ws.Cells.Item(31, 5).Value = '  +2.88%  '

# Row 32: PancakeSwap
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '2.68'
$ws.Cells.Item(32, 5).Value = '  +1.71%  '

# Row 33: NEARProtocol
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '7.02'
$ws.Cells.Item(33, 5).Value = '  -2.97%  '

# Row 34: EthereumClassic
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '28.28'
$ws.Cells.Item(34, 5).Value = '  -2.49%  '

# Row 35: Kaspa
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '0.173'
$ws.Cells.Item(35, 5).Value = '  +17.39%  '

# Row 36: Binance-PegBSC-USD
$ws.Cells.Item(36, 5).Value = '  -0.05%  '

# Row 37: RenzoRestakedETH
$ws.Cells.Item(37, 4).Value = '3.715.17'
$ws.Cells.Item(37, 5).Value = '  -1.23%  '

# Row 38: Aptos
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '8.88'
$ws.Cells.Item(38, 5).Value = '  -1.59%  '

# Row 39: Hedera
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '0.100'
$ws.Cells.Item(39, 5).Value = '  -1.21%  '

# Row 40: dogwifhat
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '3.25'
$ws.Cells.Item(40, 5).Value = '  +2.76%  '

# Row 41: Filecoin
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '5.77'
$ws.Cells.Item(41, 5).Value = '  -1.64%  '

# Row 42: Mantle -> FirstDigitalUSD
$ws.Cells.Item(42, 2).Value = 'FirstDigitalUSD'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '1.00'
$ws.Cells.Item(42, 5).Value = '  +0.06%  '

# Row 43: FirstDigitalUSD -> Mantle
$ws.Cells.Item(43, 2).Value = 'Mantle'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '0.958'
$ws.Cells.Item(43, 5).Value = '  -2.33%  '

# Row 44: USDe
$ws.Cells.Item(44, 5).Value = '  +0.01%  '

# Row 45: Monero
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '156.57'
$ws.Cells.Item(45, 5).Value = '  -0.19%  '

# Row 46: Stacks
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '1.96'
$ws.Cells.Item(46, 5).Value = '  +4.65%  '

# Row 47: ONDO
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '1.42'
$ws.Cells.Item(47, 5).Value = '  +1.72%  '

# Row 48: OKB
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '47.02'
$ws.Cells.Item(48, 5).Value = '  +1.08%  '

# Row 49: Arweave
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '42.67'
$ws.Cells.Item(49, 5).Value = '  -0.34%  '

# Row 50: TheGraph
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '0.293'
$ws.Cells.Item(50, 5).Value = '  -1.37%  '

# Row 51: Cosmos
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '8.32'
$ws.Cells.Item(51, 5).Value = '  -0.53%  '
